$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "239.43") need to be
# forced to text format first, otherwise Excel auto-converts the inline
# string into a numeric value (losing the original textual formatting).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.887.70"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "239.43"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Value = "0.6853"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.07607"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("E9").Value = "  -4.45%  "
$ws.Range("D10").Value = "23.40"
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "0.07733"
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("D12").Value = "1.825.99"
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").Value = "5.032"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "90.02"
$ws.Range("E14").Value = "  -4.42%  "
$ws.Range("D15").Value = "0.6707"
$ws.Range("E15").Value = "  -4.86%  "
$ws.Range("D16").Value = "6.404"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "0.000008260"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "28.895.79"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").Value = "242.49"
$ws.Range("E19").Value = "  -5.55%  "
$ws.Range("D20").Value = "2.091.13"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").Value = "12.60"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "7.383"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").Value = "160.59"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "8.691"
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "4.192"
$ws.Range("D31").Value = "4.145"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").Value = "1.187"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "0.05098"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").Value = "0.7528"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "1.809"
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").Value = "1.139"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "0.01827"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "1.212.18"
$ws.Range("E39").Value = "  -4.12%  "
$ws.Range("D40").Value = "2.671"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").Value = "0.9082"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "108.67"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "1.991.60"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "0.5168"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  -6.26%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "5.334"
$ws.Range("E47").Value = "  -10.20%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.429"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "62.51"
$ws.Range("E49").Value = "  -12.98%  "
$ws.Range("D50").Value = "1.719"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05852"
$ws.Range("E51").Value = "  -4.01%  "

# Remove the temporary text-format styling so these cells end up with no
# explicit style index, matching the rest of the sheet.
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
